$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1248129826832061
$ws.Range("H2").Value = 48.82120987757491
$ws.Range("I2").Value = 1.192961082138099
$ws.Range("G3").Value = 0.1128152156774726
$ws.Range("H3").Value = -3.042828234558583
$ws.Range("G4").Value = -0.03401434637986158
$ws.Range("H4").Value = -120.5104160203578
$ws.Range("G5").Value = -0.0357385225873159
$ws.Range("H5").Value = 50.26069598359501
$ws.Range("G6").Value = 0.05021042629983954
$ws.Range("H6").Value = 41.25611722415129
$ws.Range("G7").Value = 0.01782478001835862
$ws.Range("H7").Value = -12.17199474912758
$ws.Range("G8").Value = -0.1362091962515835
$ws.Range("H8").Value = 3.464757877876811
$ws.Range("G9").Value = -0.1162330924337188
$ws.Range("H9").Value = 15.08829636322568
$ws.Range("G10").Value = -0.08234249921738504
$ws.Range("H10").Value = 23.77608805466782
$ws.Range("G11").Value = -0.09095939367564763
$ws.Range("H11").Value = -36.69704815614066
$ws.Range("G12").Value = -0.3741618449121005
$ws.Range("H12").Value = 9.883150738201937
$ws.Range("G13").Value = -0.4062934356010162
$ws.Range("H13").Value = 9.49120305582502
$ws.Range("G14").Value = -0.03086123641758959
$ws.Range("H14").Value = 39.17253759619449
$ws.Range("G15").Value = -0.001357785827895729
$ws.Range("H15").Value = 98.35800915976466
$ws.Range("G16").Value = 0.1208821602501296
$ws.Range("H16").Value = -17.21686285131582
$ws.Range("G17").Value = 0.1443809287207572
$ws.Range("H17").Value = 17.74521274660619
$ws.Range("G18").Value = 0.1363484186363655
$ws.Range("H18").Value = -1.523042718420635
$ws.Range("G19").Value = 0.1433937589176967
$ws.Range("H19").Value = 50.29729732916839
$ws.Range("G20").Value = 0.01360042559572368
$ws.Range("H20").Value = -46.94391070477423
$ws.Range("G21").Value = 0.03006223478057242
$ws.Range("H21").Value = -60.00613330736326
$ws.Range("G24").Value = 0.09743356845847011
$ws.Range("H24").Value = -2.989887446746384
$ws.Range("G25").Value = 0.1520578091251967
$ws.Range("H25").Value = 0.3373717024817708
$ws.Range("G26").Value = 0.06510234227117315
$ws.Range("H26").Value = -17.72185825869026
$ws.Range("G27").Value = 0.06249001417556976
$ws.Range("H27").Value = -37.41775278244013
$ws.Range("G28").Value = -0.2392594294414418
$ws.Range("H28").Value = -12.26555221595133
$ws.Range("G29").Value = -0.2016837271915152
$ws.Range("H29").Value = 1.764652143604634
$ws.Range("G30").Value = 0.06550112217549428
$ws.Range("H30").Value = 48.42212578568261
$ws.Range("G31").Value = 0.02752753816999139
$ws.Range("H31").Value = 4.535834291031299
$ws.Range("G32").Value = 0.09832908843531273
$ws.Range("H32").Value = 3.547347960879966
$ws.Range("G33").Value = 0.1607910243468573
$ws.Range("H33").Value = 54.65897098069754
$ws.Range("G34").Value = 0.02765813665890893
$ws.Range("H34").Value = -40.42749944416458
$ws.Range("G35").Value = 0.01946275077790345
$ws.Range("H35").Value = 156.8822018113669
$ws.Range("G36").Value = 0.04525559492812234
$ws.Range("H36").Value = -21.61666523347264
$ws.Range("G37").Value = 0.06540431242063113
$ws.Range("H37").Value = -6.998600512074622
$ws.Range("G38").Value = 0.01182934197795375
$ws.Range("H38").Value = -77.41906408752708
$ws.Range("G39").Value = 0.03424010456724188
$ws.Range("H39").Value = 65.11337611887383
$ws.Range("G40").Value = 0.005189021683586412
$ws.Range("H40").Value = 161.1120961917191
$ws.Range("G41").Value = 0.04550105144361731
$ws.Range("H41").Value = 28.6945336672844
$ws.Range("G42").Value = 0.1399546251547427
$ws.Range("H42").Value = 4.688967783887578
$ws.Range("G43").Value = 0.1422085549869157
$ws.Range("H43").Value = -4.542855431958436
$ws.Range("G44").Value = -0.007153036006435855
$ws.Range("H44").Value = 15.96015101852601
$ws.Range("G45").Value = 0.008495153769514397
$ws.Range("H45").Value = 177.3702930160298
$ws.Range("G46").Value = -0.002422132313178876
$ws.Range("H46").Value = 26.44705296175452
$ws.Range("G47").Value = 0.000251761457700876
$ws.Range("H47").Value = 102.7133908569362
$ws.Range("G48").Value = 0.04620871750503409
$ws.Range("H48").Value = -8.082800472770577
$ws.Range("G49").Value = 0.06375211059374096
$ws.Range("H49").Value = -3.500613961721503
$ws.Range("G50").Value = 0.1290134155123612
$ws.Range("H50").Value = -19.99549090749029
$ws.Range("G51").Value = 0.1587331114160908
$ws.Range("H51").Value = -7.235982227304513
$ws.Range("G52").Value = -0.1634145080841927
$ws.Range("H52").Value = -1.859682004574498
$ws.Range("G53").Value = -0.1585205914821631
$ws.Range("H53").Value = -25.75559859960287
$ws.Range("G54").Value = 0.1129299203335908
$ws.Range("H54").Value = 20.49406131106061
$ws.Range("G55").Value = 0.0919386737899725
$ws.Range("H55").Value = -18.69643903214667
$ws.Range("G56").Value = -0.01697517127572675
$ws.Range("H56").Value = -132.5120196623836
$ws.Range("G57").Value = -0.01399544675128893
$ws.Range("H57").Value = 38.79736169433686
$ws.Range("G58").Value = 0.03136438773211132
$ws.Range("H58").Value = -44.37752182490204
$ws.Range("G59").Value = 0.07474663768395225
$ws.Range("H59").Value = 4.075069243325852
$ws.Range("G60").Value = 0.05278052175225122
$ws.Range("H60").Value = -24.57118805758202
$ws.Range("G61").Value = 0.07197855302671317
$ws.Range("H61").Value = 51.44839320337438
$ws.Range("G62").Value = 0.05945020004304512
$ws.Range("H62").Value = -18.51708720277197
$ws.Range("G63").Value = 0.0661072257244885
$ws.Range("H63").Value = 1.091425191106748
$ws.Range("G64").Value = -0.03198310920447313
$ws.Range("H64").Value = 22.77503985443016
$ws.Range("G65").Value = -0.01380218701311709
$ws.Range("H65").Value = 72.01713403268
$ws.Range("G66").Value = 0.04908219371731174
$ws.Range("H66").Value = 159.2281476093465
$ws.Range("G67").Value = 0.00761273890300986
$ws.Range("H67").Value = -70.88806067691603
$ws.Range("G68").Value = -0.006687354134776244
$ws.Range("H68").Value = -1273.315996255327
$ws.Range("G69").Value = 0.02249636686908364
$ws.Range("H69").Value = 274.2242283320046
$ws.Range("G70").Value = -0.0457433628228579
$ws.Range("H70").Value = -66.6542163488776
$ws.Range("G71").Value = -0.06137741277376518
$ws.Range("H71").Value = -11.40903829152097
$ws.Range("G72").Value = -0.1740632958065455
$ws.Range("H72").Value = -17.34539908155569
$ws.Range("G73").Value = -0.1705668012824217
$ws.Range("H73").Value = -17.80302125538277
$ws.Range("G74").Value = 0.1495384657514329
$ws.Range("H74").Value = 18.63320436351996
$ws.Range("G75").Value = 0.1479259838466666
$ws.Range("H75").Value = 9.435392522727865
$ws.Range("G76").Value = -0.04372318294825396
$ws.Range("H76").Value = -26.96222987779667
$ws.Range("G77").Value = -0.03329774771581427
$ws.Range("H77").Value = 27.90667923030474
$ws.Range("G78").Value = 0.08852303452429428
$ws.Range("H78").Value = -3.95838425241352
$ws.Range("G79").Value = 0.0882313592585698
$ws.Range("H79").Value = -8.572736256026831
$ws.Range("G80").Value = -0.1927557997519705
$ws.Range("H80").Value = -18.68259795429503
$ws.Range("G81").Value = -0.1836457056487156
$ws.Range("H81").Value = 15.15499471321745
$ws.Range("G82").Value = 0.1515716080431263
$ws.Range("H82").Value = 9.243728221890441
$ws.Range("G83").Value = 0.2034107166727244
$ws.Range("H83").Value = 23.56409739613777
$ws.Range("G84").Value = 0.05423850344719331
$ws.Range("H84").Value = 287.5555315228483
$ws.Range("G85").Value = 0.05490533035735638
$ws.Range("H85").Value = 142.5234543538503
